# "user can now include units in the inputs"
# The input-parameter headers on the "values" sheet get their unit of
# measure appended in parentheses so the spreadsheet is self-documenting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("values")
$ws.Activate()

$ws.Range("D1").Value = "w_p_input(lb)"
$ws.Range("K1").Value = "capital_a_input(in)"
$ws.Range("L1").Value = "capital_b_input(in)"
$ws.Range("M1").Value = "a_input(in)"
$ws.Range("N1").Value = "b_input(in)"
$ws.Range("O1").Value = "capital_h_input(ft)"

# Scroll the view over a couple of columns and land the selection on I25,
# matching where the editor ended up after making the change.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("I25").Select()
